$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as text, matching
# the workbook's existing inline-string cell type.
function Set-CellText {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

$ws.Range("D2").Value = "62.815.86"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").Value = "2.462.66"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("E4").Value = "  -0.04%  "

Set-CellText $ws.Range("D5") "574.72"
$ws.Range("E5").Value = "  -0.65%  "

Set-CellText $ws.Range("D6") "145.78"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").Value = "2.461.96"
$ws.Range("E9").Value = "  +0.66%  "

$ws.Range("E10").Value = "  +1.30%  "

Set-CellText $ws.Range("D11") "0.163"
$ws.Range("E11").Value = "  +1.48%  "

Set-CellText $ws.Range("D12") "5.27"
$ws.Range("E12").Value = "  +0.78%  "

$ws.Range("E13").Value = "  +0.62%  "

Set-CellText $ws.Range("D14") "28.98"
$ws.Range("E14").Value = "  +2.09%  "

$ws.Range("E15").Value = "  -0.46%  "

$ws.Range("D16").Value = "2.910.18"
$ws.Range("E16").Value = "  +0.61%  "

$ws.Range("D17").Value = "62.775.20"
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("D18").Value = "2.457.06"
$ws.Range("E18").Value = "  +0.64%  "

Set-CellText $ws.Range("D19") "8.00"
$ws.Range("E19").Value = "  +2.47%  "

Set-CellText $ws.Range("D20") "11.00"
$ws.Range("E20").Value = "  +0.64%  "

Set-CellText $ws.Range("D21") "326.93"
$ws.Range("E21").Value = "  +0.24%  "

$ws.Range("B22").Value = "SuiNetwork"
$ws.Range("C22").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-CellText $ws.Range("D22") "2.22"
$ws.Range("E22").Value = "  +10.11%  "

$ws.Range("B23").Value = "Polkadot"
$ws.Range("C23").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-CellText $ws.Range("D23") "4.13"
$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("E24").Value = "  +0.02%  "

Set-CellText $ws.Range("D25") "10.27"
$ws.Range("E25").Value = "  +20.06%  "

Set-CellText $ws.Range("D26") "65.73"
$ws.Range("E26").Value = "  +0.63%  "

Set-CellText $ws.Range("D27") "653.96"
$ws.Range("E27").Value = "  +1.10%  "

$ws.Range("D28").Value = "0.0₃0980"
$ws.Range("E28").Value = "  +0.27%  "

$ws.Range("D29").Value = "2.583.05"

$ws.Range("E30").Value = "  -12.44%  "

$ws.Range("E31").Value = "  +2.46%  "

Set-CellText $ws.Range("D32") "7.98"
$ws.Range("E32").Value = "  -2.40%  "

$ws.Range("E33").Value = "  -1.00%  "

Set-CellText $ws.Range("D34") "0.133"
$ws.Range("E34").Value = "  -3.74%  "

Set-CellText $ws.Range("D35") "0.999"
$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("E36").Value = "  +3.01%  "

Set-CellText $ws.Range("D38") "0.369"
$ws.Range("E38").Value = "  -1.01%  "

$ws.Range("E39").Value = "  +0.69%  "

Set-CellText $ws.Range("D40") "151.03"
$ws.Range("E40").Value = "  -1.23%  "

Set-CellText $ws.Range("D41") "5.37"
$ws.Range("E41").Value = "  -1.79%  "

Set-CellText $ws.Range("D42") "2.76"
$ws.Range("E42").Value = "  +2.01%  "

$ws.Range("E43").Value = "  -1.53%  "

$ws.Range("E44").Value = "  -71.83%  "

$ws.Range("E45").Value = "  -0.07%  "

Set-CellText $ws.Range("D46") "154.32"
$ws.Range("E46").Value = "  +6.95%  "

Set-CellText $ws.Range("D47") "15.22"
$ws.Range("E47").Value = "  +1.37%  "

Set-CellText $ws.Range("D48") "3.58"
$ws.Range("E48").Value = "  -0.17%  "

Set-CellText $ws.Range("D49") "20.29"
$ws.Range("E49").Value = "  -1.36%  "

Set-CellText $ws.Range("D50") "0.607"
$ws.Range("E50").Value = "  +0.52%  "

Set-CellText $ws.Range("D51") "0.0511"
$ws.Range("E51").Value = "  -0.25%  "
